$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Project" score for assignment #1 (row 8)
$ws.Range("E8").Value = 0.925

# Update the "HW" score for assignment #2 (row 9)
$ws.Range("D9").Value = 0.7

# Reflect the final selection made by the user (clicked on D10 afterward)
$ws.Range("D10").Select()
